$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C2:C15) from 2023-09-09 (45178)
# to 2023-09-10 (45179), keeping the existing date formatting/style.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
